# repull data, push all data, mean calculation
# Update the dSF column (F) values for the rows whose underlying data was
# repulled, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    6  = -3
    8  = -1
    9  = -4
    14 = -1
    15 = -1
    21 = -9
    25 = -6
    26 = -3
    28 = -1
    31 = -2
    33 = -2
    40 = 3
    43 = -2
    47 = -1
    51 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
